# Updates the cryptos price list (Price / Volume(1h) columns) with the
# latest scrape values, matching the GitHub Actions "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "55.108.95"
$ws.Range("E2").Value = "  -3.69%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.938.95"
$ws.Range("E3").Value = "  -6.14%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "483.68"
$ws.Range("E5").Value = "  -6.50%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.01"
$ws.Range("E6").Value = "  -3.04%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.935.66"
$ws.Range("E8").Value = "  -6.28%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.413"
$ws.Range("E9").Value = "  -7.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.94"
$ws.Range("E10").Value = "  -3.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0991"
$ws.Range("E11").Value = "  -7.72%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.344"
$ws.Range("E12").Value = "  -9.61%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.127"
$ws.Range("E13").Value = "  +0.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.449.95"
$ws.Range("E14").Value = "  -6.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.24"
$ws.Range("E15").Value = "  -3.53%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "55.154.05"
$ws.Range("E16").Value = "  -3.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.945.70"
$ws.Range("E17").Value = "  -6.20%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000137"
$ws.Range("E18").Value = "  -7.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.54"
$ws.Range("E19").Value = "  -3.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.84"
$ws.Range("E20").Value = "  -7.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.33"
$ws.Range("E21").Value = "  -7.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "310.00"
$ws.Range("E22").Value = "  -9.28%  "

$ws.Range("E23").Value = "  +0.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.457"
$ws.Range("E24").Value = "  -8.81%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "59.58"
$ws.Range("E25").Value = "  -13.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("E27").Value = "  -3.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0827"
$ws.Range("E29").Value = "  -10.86%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.46"
$ws.Range("E30").Value = "  -2.86%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.41"
$ws.Range("E31").Value = "  -6.01%  "

$ws.Range("E32").Value = "  -3.79%  "

$ws.Range("E33").Value = "  -9.44%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.09"
$ws.Range("E34").Value = "  -11.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "147.47"
$ws.Range("E35").Value = "  -5.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.31"
$ws.Range("E36").Value = "  -10.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.58"
$ws.Range("E37").Value = "  -9.28%  "

$ws.Range("E38").Value = "  -8.48%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "23.21"
$ws.Range("E39").Value = "  -9.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0637"
$ws.Range("E40").Value = "  -6.65%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.977.61"
$ws.Range("E41").Value = "  -5.95%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "35.73"
$ws.Range("E43").Value = "  -11.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.625"
$ws.Range("E44").Value = "  -9.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.974"
$ws.Range("E45").Value = "  -7.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.36"
$ws.Range("E46").Value = "  -6.64%  "

$ws.Range("E47").Value = "  -9.93%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.109.06"
$ws.Range("E48").Value = "  -5.56%  "

$ws.Range("E49").Value = "  -1.41%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.67"
$ws.Range("E50").Value = "  -6.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.50"
$ws.Range("E51").Value = "  -10.14%  "
